$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Report")
$ws.Range("AW2").Value = 203.917326
$ws.Range("AW3").Value = 47.033322
$ws.Range("AW4").Value = 34.791366
$ws.Range("AW5").Value = 158.862986
$ws.Range("AW6").Value = 161.649537
$ws.Range("AQ7").Value = 25.90515
$ws.Range("AW8").Value = 196.687037
$ws.Range("AW9").Value = 145.727674
$ws.Range("AW10").Value = 129.810961
$ws.Range("AW11").Value = 45.707604
$ws.Range("AW12").Value = 130.713333
$ws.Range("AW13").Value = 208.77103
$ws.Range("AW14").Value = 194.586181
$ws.Range("AW15").Value = 163.714711
$ws.Range("AW16").Value = 166.012014
$ws.Range("AW17").Value = 45.661829
$ws.Range("AW18").Value = 27.89809
$ws.Range("AW19").Value = 133.598843
$ws.Range("AW20").Value = 125.687604
$ws.Range("AW21").Value = 47.03375
$ws.Range("AQ22").Value = 25.904942
$ws.Range("AQ23").Value = 25.769306
$ws.Range("AW24").Value = 158.862454
$ws.Range("AW25").Value = 45.806667
$ws.Range("AT26").Value = 12.898044
$ws.Range("AW27").Value = 194.59169
$ws.Range("AW28").Value = 115.675174
$ws.Range("AW29").Value = 124.768657
$ws.Range("AW30").Value = 54.908704
$ws.Range("AW31").Value = 196.681944
$ws.Range("AW32").Value = 53.690104
$ws.Range("AW33").Value = 132.91066
$ws.Range("AW34").Value = 41.808449
$ws.Range("AW35").Value = 160.956458
$ws.Range("AW36").Value = 168.742199
$ws.Range("AW37").Value = 55.974097
$ws.Range("AK38").Value = 61.754688
$ws.Range("AW39").Value = 34.788345
$ws.Range("AK40").Value = 61.754583
$ws.Range("AQ41").Value = 27.17662
$ws.Range("AN42").Value = 33.731296
$ws.Range("AW43").Value = 117.993461
$ws.Range("AW44").Value = 132.910046
$ws.Range("AW45").Value = 84.02681699999999
$ws.Range("AW46").Value = 84.02560200000001
$ws.Range("AK47").Value = 83.80420100000001
$ws.Range("AW48").Value = 160.957442
$ws.Range("AW49").Value = 161.650949
$ws.Range("AW50").Value = 68.793565
$ws.Range("AW51").Value = 174.866435
$ws.Range("AW52").Value = 152.62397
$ws.Range("AW53").Value = 199.728623
$ws.Range("AW54").Value = 189.777442
$ws.Range("AW55").Value = 19.879988
$ws.Range("AW56").Value = 161.650116
$ws.Range("AW57").Value = 55.978206
$ws.Range("AW58").Value = 35.597847
$ws.Range("AW59").Value = 97.65169
$ws.Range("AW60").Value = 164.905799
$ws.Range("AW61").Value = 125.706944
$ws.Range("AW62").Value = 45.706817
$ws.Range("AW63").Value = 201.021759
$ws.Range("AW64").Value = 203.914653
$ws.Range("AW65").Value = 161.665891
$ws.Range("AW66").Value = 118.767419
$ws.Range("AW67").Value = 54.804016
$ws.Range("AW68").Value = 61.799977
$ws.Range("AW69").Value = 18.62522
$ws.Range("AQ70").Value = 24.928241
$ws.Range("AW71").Value = 13.927674
$ws.Range("AW72").Value = 175.737326
$ws.Range("AW73").Value = 112.763542
$ws.Range("AW74").Value = 47.690532
$ws.Range("AK75").Value = 18.781968
$ws.Range("AW76").Value = 116.972188
$ws.Range("AW77").Value = 34.778692
$ws.Range("AW78").Value = 189.780833
$ws.Range("AW79").Value = 84.026597
$ws.Range("AW80").Value = 68.941956
$ws.Range("AW81").Value = 63.636713
$ws.Range("AW82").Value = 53.688762
$ws.Range("AW83").Value = 168.806157
$ws.Range("AW84").Value = 199.680903
$ws.Range("AW85").Value = 147.650567
$ws.Range("AW86").Value = 116.970475
$ws.Range("AW87").Value = 97.644537
$ws.Range("AW88").Value = 19.889352
$ws.Range("AW89").Value = 153.653079
$ws.Range("AW90").Value = 194.585799
$ws.Range("AW91").Value = 104.026227
$ws.Range("AW92").Value = 161.64934
$ws.Range("AT93").Value = 10.837847
$ws.Range("AW94").Value = 201.022153
$ws.Range("AW95").Value = 112.636817
$ws.Range("AW96").Value = 34.791921
$ws.Range("AW97").Value = 130.705995
$ws.Range("AW98").Value = 122.781852
$ws.Range("AW99").Value = 196.771759
$ws.Range("AW100").Value = 152.76875
$ws.Range("AW101").Value = 92.0275
$ws.Range("AW102").Value = 33.806817
$ws.Range("AW103").Value = 151.205625
$ws.Range("AW104").Value = 159.849144
$ws.Range("AW105").Value = 112.692766
$ws.Range("AW106").Value = 61.781458
$ws.Range("AW107").Value = 34.777373
$ws.Range("AQ108").Value = 32.815567
$ws.Range("AT109").Value = 10.837419
$ws.Range("AW110").Value = 203.917616
$ws.Range("AW111").Value = 189.728495
$ws.Range("AW112").Value = 199.679039
$ws.Range("AW113").Value = 92.025903
$ws.Range("AW114").Value = 63.62735
$ws.Range("AW115").Value = 45.806181
$ws.Range("AW116").Value = 178.647708
$ws.Range("AW117").Value = 152.788762
$ws.Range("AW118").Value = 130.659201
$ws.Range("AW119").Value = 116.968322
$ws.Range("AK120").Value = 61.754699
$ws.Range("AW121").Value = 194.586412
$ws.Range("AW122").Value = 34.778796
$ws.Range("AW123").Value = 18.625926
$ws.Range("AK124").Value = 18.781956
$ws.Range("AW125").Value = 189.681007
$ws.Range("AW126").Value = 19.881435
$ws.Range("AK127").Value = 61.754583
$ws.Range("AW128").Value = 103.888576
$ws.Range("AW129").Value = 98.03998799999999
$ws.Range("AW130").Value = 207.73103
$ws.Range("AW131").Value = 133.571389
$ws.Range("AW132").Value = 92.02535899999999
$ws.Range("AQ133").Value = 32.703785
$ws.Range("AW134").Value = 27.897859
$ws.Range("AW135").Value = 189.731701
$ws.Range("AW136").Value = 194.588067
$ws.Range("AW137").Value = 130.768692
$ws.Range("AK138").Value = 61.754618
$ws.Range("AW139").Value = 164.90537
$ws.Range("AW140").Value = 161.649838
$ws.Range("AW141").Value = 124.769688
$ws.Range("AQ142").Value = 15.022917
$ws.Range("AW143").Value = 96.90186300000001
$ws.Range("AW144").Value = 96.701921
$ws.Range("AW145").Value = 19.888981
$ws.Range("AW146").Value = 118.600058
$ws.Range("AW147").Value = 108.665289
$ws.Range("AW148").Value = 138.735463
$ws.Range("AW149").Value = 92.024664
$ws.Range("AW150").Value = 35.598507
$ws.Range("AQ151").Value = 25.860463
$ws.Range("AW152").Value = 160.957245
$ws.Range("AW153").Value = 116.973935
$ws.Range("AW154").Value = 40.621458
$ws.Range("AW155").Value = 150.980058
$ws.Range("AW156").Value = 66.726505
$ws.Range("AW157").Value = 34.775718
$ws.Range("AW158").Value = 199.68184
$ws.Range("AW159").Value = 176.032141
$ws.Range("AW160").Value = 115.674965
$ws.Range("AT161").Value = 11.571296
$ws.Range("AW162").Value = 150.985
$ws.Range("AW163").Value = 203.918009
$ws.Range("AW164").Value = 125.686076
$ws.Range("AW165").Value = 196.681389
$ws.Range("AW166").Value = 187.788229
$ws.Range("AW167").Value = 125.706829
$ws.Range("AW168").Value = 54.803854
$ws.Range("AK169").Value = 83.80419000000001
$ws.Range("AW170").Value = 196.68103
$ws.Range("AW171").Value = 152.784572
$ws.Range("AW172").Value = 98.039271
$ws.Range("AW173").Value = 130.598981
$ws.Range("AW174").Value = 41.752326
$ws.Range("AW175").Value = 203.925579
$ws.Range("AW176").Value = 138.74765
$ws.Range("AW177").Value = 116.971447
$ws.Range("AW178").Value = 98.042813
$ws.Range("AW179").Value = 152.788391
$ws.Range("AK180").Value = 61.754699
$ws.Range("AW181").Value = 19.880752
$ws.Range("AW182").Value = 196.772037
$ws.Range("AW183").Value = 210.99287
$ws.Range("AW184").Value = 116.977292
$ws.Range("AW185").Value = 90.76468800000001
$ws.Range("AW186").Value = 184.008576
$ws.Range("AW187").Value = 194.586979
$ws.Range("AW188").Value = 122.690532
$ws.Range("AW189").Value = 178.644549
$ws.Range("AW190").Value = 126.61397
$ws.Range("AW191").Value = 199.655174
$ws.Range("AW192").Value = 196.771493
$ws.Range("AW193").Value = 195.693715
$ws.Range("AW194").Value = 207.676296
$ws.Range("AW195").Value = 209.81669
$ws.Range("AW196").Value = 199.729873
$ws.Range("AW197").Value = 161.651563
$ws.Range("AW198").Value = 66.99457200000001
$ws.Range("AW199").Value = 161.669537
$ws.Range("AK200").Value = 83.80420100000001
$ws.Range("AW201").Value = 208.772014
$ws.Range("AW202").Value = 199.670938
$ws.Range("AW203").Value = 208.767882
$ws.Range("AW204").Value = 109.722188
$ws.Range("AW205").Value = 41.776481
$ws.Range("AW206").Value = 34.79463
$ws.Range("AQ207").Value = 25.905058
$ws.Range("AQ208").Value = 25.905255
$ws.Range("AW209").Value = 116.976644
$ws.Range("AW210").Value = 168.742118
$ws.Range("AW211").Value = 185.983426
$ws.Range("AW212").Value = 125.753032
$ws.Range("AW213").Value = 117.993229
$ws.Range("AW214").Value = 98.041123
$ws.Range("AW215").Value = 48.969988
$ws.Range("AW216").Value = 160.957778
$ws.Range("AW217").Value = 116.836238
$ws.Range("AW218").Value = 19.881979
$ws.Range("AW219").Value = 35.599005
$ws.Range("AW220").Value = 150.983773
$ws.Range("AW221").Value = 196.656273
$ws.Range("AW222").Value = 161.650347
$ws.Range("AW223").Value = 125.684421
$ws.Range("AW224").Value = 68.93900499999999
$ws.Range("AW225").Value = 97.743889
$ws.Range("AW226").Value = 92.023785
$ws.Range("AK227").Value = 61.754595
$ws.Range("AW228").Value = 161.650752
$ws.Range("AW229").Value = 132.911285
$ws.Range("AN230").Value = 18.802546
$ws.Range("AW231").Value = 181.803947
$ws.Range("AQ232").Value = 32.779861
$ws.Range("AT233").Value = 11.606771
$ws.Range("AW234").Value = 160.918356
$ws.Range("AW235").Value = 92.03353
$ws.Range("AW236").Value = 92.03137700000001
$ws.Range("AW237").Value = 66.615926
$ws.Range("AW238").Value = 158.862616
$ws.Range("AW239").Value = 83.02209499999999
$ws.Range("AW240").Value = 68.79431700000001
$ws.Range("AW241").Value = 17.791007
$ws.Range("AW242").Value = 189.781053
$ws.Range("AW243").Value = 199.729132
$ws.Range("AW244").Value = 178.644896
$ws.Range("AW245").Value = 54.80375
$ws.Range("AW246").Value = 68.79512699999999
$ws.Range("AW247").Value = 34.778519
$ws.Range("AW248").Value = 178.64544
$ws.Range("AW249").Value = 98.035521
$ws.Range("AW250").Value = 49.782361
$ws.Range("AW251").Value = 188.674248
$ws.Range("AW252").Value = 103.889769
$ws.Range("AW253").Value = 34.795289
$ws.Range("AW254").Value = 55.976377
$ws.Range("AW255").Value = 174.866238
$ws.Range("AW256").Value = 54.909433
$ws.Range("AW257").Value = 45.70544
$ws.Range("AW258").Value = 28.738148
$ws.Range("AW259").Value = 84.025521
$ws.Range("AW260").Value = 176.032512
$ws.Range("AW261").Value = 187.789178
$ws.Range("AK262").Value = 18.781956
$ws.Range("AW263").Value = 150.984838
$ws.Range("AW264").Value = 130.712905
$ws.Range("AW265").Value = 72.040324
$ws.Range("AW266").Value = 49.776019
$ws.Range("AW267").Value = 174.994907
$ws.Range("AK268").Value = 83.804213
$ws.Range("AW269").Value = 188.994109
$ws.Range("AW270").Value = 133.602002
$ws.Range("AW271").Value = 194.585012
$ws.Range("AW272").Value = 116.977859
$ws.Range("AW273").Value = 161.652685
$ws.Range("AW274").Value = 34.646748
$ws.Range("AW275").Value = 185.983889
$ws.Range("AW276").Value = 116.970729
$ws.Range("AW277").Value = 91.78219900000001
$ws.Range("AW278").Value = 13.927917
$ws.Range("AW279").Value = 126.743669
$ws.Range("AW280").Value = 201.023611
$ws.Range("AW281").Value = 90.69582200000001
$ws.Range("AW282").Value = 45.80603
$ws.Range("AW283").Value = 19.893495
$ws.Range("AT284").Value = 10.837535
$ws.Range("AW285").Value = 90.69438700000001
$ws.Range("AW286").Value = 90.696157
$ws.Range("AW287").Value = 194.585394
$ws.Range("AW288").Value = 147.647083
$ws.Range("AW289").Value = 66.72745399999999
$ws.Range("AW290").Value = 147.649873
$ws.Range("AW291").Value = 160.955116
$ws.Range("AW292").Value = 203.918461
$ws.Range("AW293").Value = 98.03466400000001
$ws.Range("AW294").Value = 61.744745
$ws.Range("AW295").Value = 49.78
$ws.Range("AW296").Value = 182.692512
$ws.Range("AW297").Value = 137.751528
$ws.Range("AW298").Value = 19.895023
$ws.Range("AQ299").Value = 27.176493
$ws.Range("AW300").Value = 161.652014
$ws.Range("AW301").Value = 130.6589
$ws.Range("AW302").Value = 76.669618
$ws.Range("AW303").Value = 55.972593
$ws.Range("AW304").Value = 45.705868
$ws.Range("AW305").Value = 140.676493
$ws.Range("AN306").Value = 20.613241
$ws.Range("AW307").Value = 132.907928
$ws.Range("AW308").Value = 66.726968
$ws.Range("AW309").Value = 208.77684
$ws.Range("AW310").Value = 118.97581
$ws.Range("AW311").Value = 98.042025
$ws.Range("AW312").Value = 47.690046
$ws.Range("AW313").Value = 166.015185
$ws.Range("AQ314").Value = 32.665313
$ws.Range("AK315").Value = 61.75478
$ws.Range("AW316").Value = 28.738762
$ws.Range("AW317").Value = 122.735567
$ws.Range("AW318").Value = 101.674086
$ws.Range("AW319").Value = 66.72807899999999
$ws.Range("AW320").Value = 47.691273
$ws.Range("AW321").Value = 174.995127
$ws.Range("AW322").Value = 103.889317
$ws.Range("AW323").Value = 92.023021
$ws.Range("AW324").Value = 98.09553200000001
$ws.Range("AW325").Value = 33.810278
$ws.Range("AQ326").Value = 15.023067
$ws.Range("AW327").Value = 199.679745
$ws.Range("AW328").Value = 201.021921
$ws.Range("AK329").Value = 61.754745
$ws.Range("AW330").Value = 14.809213
$ws.Range("AW331").Value = 207.67059
$ws.Range("AW332").Value = 98.04324099999999
$ws.Range("AW333").Value = 17.790926
$ws.Range("AW334").Value = 178.645961
$ws.Range("AW335").Value = 68.94050900000001
$ws.Range("AW336").Value = 66.975544
$ws.Range("AW337").Value = 194.584745
$ws.Range("AW338").Value = 166.0125
$ws.Range("AW339").Value = 84.612697
$ws.Range("AW340").Value = 168.805058
$ws.Range("AW341").Value = 116.972083
$ws.Range("AW342").Value = 203.810498
$ws.Range("AW343").Value = 199.716354
$ws.Range("AW344").Value = 47.696667
$ws.Range("AW345").Value = 201.02272
$ws.Range("AW346").Value = 178.646481
$ws.Range("AW347").Value = 208.768241
$ws.Range("AW348").Value = 56.704618
$ws.Range("AW349").Value = 194.583275
$ws.Range("AW350").Value = 161.649155
$ws.Range("AW351").Value = 152.781782
$ws.Range("AW352").Value = 68.939826
$ws.Range("AW353").Value = 19.893947
$ws.Range("AW354").Value = 13.926759
$ws.Range("AW355").Value = 115.675313
$ws.Range("AW356").Value = 115.67485
$ws.Range("AW357").Value = 19.892373
$ws.Range("AW358").Value = 34.792894
$ws.Range("AW359").Value = 47.690255
